$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.304.51'
$ws.Range('E2').Value = '  -7.02%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.718.82'
$ws.Range('E3').Value = '  -6.77%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.04'
$ws.Range('E5').Value = '  -6.30%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.65'
$ws.Range('E6').Value = '  +6.31%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.713.08'
$ws.Range('E7').Value = '  -6.72%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.637'
$ws.Range('E8').Value = '  -6.84%  '

# Row 9
$ws.Range('E9').Value = '  +0.03%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.720'
$ws.Range('E10').Value = '  -4.48%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.168'
$ws.Range('E11').Value = '  -10.62%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.28'
$ws.Range('E12').Value = '  -5.40%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000304'
$ws.Range('E13').Value = '  -11.13%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.75'
$ws.Range('E14').Value = '  -2.96%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.314.49'
$ws.Range('E15').Value = '  -6.73%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.758.24'
$ws.Range('E16').Value = '  -6.04%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.54'
$ws.Range('E17').Value = '  -4.76%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.128'
$ws.Range('E18').Value = '  -2.79%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.11'
$ws.Range('E19').Value = '  -7.29%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.15'
$ws.Range('E20').Value = '  -7.87%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.169.77'
$ws.Range('E21').Value = '  -7.04%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '412.77'
$ws.Range('E22').Value = '  -6.40%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.61'
$ws.Range('E23').Value = '  -6.01%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.35'
$ws.Range('E24').Value = '  -6.80%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.10'
$ws.Range('E25').Value = '  -7.84%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.96'
$ws.Range('E26').Value = '  -8.42%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.89'
$ws.Range('E27').Value = '  -1.89%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.85'
$ws.Range('E28').Value = '  -5.10%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.94'
$ws.Range('E29').Value = '  +0.06%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.70'
$ws.Range('E30').Value = '  -7.15%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.28'
$ws.Range('E31').Value = '  +6.52%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.27'
$ws.Range('E32').Value = '  -7.93%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.87'
$ws.Range('E33').Value = '  -5.91%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '45.27'
$ws.Range('E34').Value = '  -5.04%  '

# Row 35
$ws.Range('E35').Value = '  -8.42%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '66.27'
$ws.Range('E36').Value = '  -7.05%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0942'
$ws.Range('E37').Value = '  -5.24%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '600.43'
$ws.Range('E38').Value = '  -5.78%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.407'
$ws.Range('E39').Value = '  -5.72%  '

# Row 40
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.12%  '

# Row 41
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.09%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.27'
$ws.Range('E42').Value = '  +14.03%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.137'
$ws.Range('E43').Value = '  -6.61%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.09'
$ws.Range('E44').Value = '  -9.27%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0445'
$ws.Range('E45').Value = '  -8.38%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.62'
$ws.Range('E46').Value = '  -13.22%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.62'
$ws.Range('E47').Value = '  +0.40%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.136'
$ws.Range('E48').Value = '  -8.43%  '

# Row 49
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.762.34'
$ws.Range('E49').Value = '  -3.47%  '

# Row 50
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.16'
$ws.Range('E50').Value = '  -7.76%  '

# Row 51
$ws.Range('E51').Value = '  -4.74%  '
